$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated win/loss/transition probabilities from games pulled March 7
# (team_specific_matrix/Creighton_B.xlsx)

# Row 2
$ws.Range("B2").Value = 0.1676646706586826
$ws.Range("C2").Value = 0.6526946107784432
$ws.Range("P2").Value = 0.08982035928143713
$ws.Range("S2").Value = 0.08982035928143713

# Row 3
$ws.Range("C3").Value = 0.05263157894736842
$ws.Range("J3").Value = 0.01754385964912281
$ws.Range("P3").Value = 0.8070175438596491
$ws.Range("S3").Value = 0.1228070175438596

# Row 4
$ws.Range("P4").Value = 0.6774193548387096
$ws.Range("S4").Value = 0.3225806451612903

# Row 6
$ws.Range("B6").Value = 0.03804347826086957
$ws.Range("D6").Value = 0.0108695652173913
$ws.Range("F6").Value = 0.04347826086956522
$ws.Range("J6").Value = 0.2119565217391304
$ws.Range("O6").Value = 0.03260869565217391
$ws.Range("Q6").Value = 0.2173913043478261
$ws.Range("R6").Value = 0.06521739130434782
$ws.Range("S6").Value = 0.3804347826086957

# Row 7
$ws.Range("B7").Value = 0.06629834254143646
$ws.Range("D7").Value = 0.01104972375690608
$ws.Range("F7").Value = 0.0718232044198895
$ws.Range("J7").Value = 0.1270718232044199
$ws.Range("O7").Value = 0.02762430939226519
$ws.Range("Q7").Value = 0.2541436464088398
$ws.Range("R7").Value = 0.04972375690607735
$ws.Range("S7").Value = 0.3922651933701657

# Row 8
$ws.Range("B8").Value = 0.0396039603960396
$ws.Range("D8").Value = 0.01782178217821782
$ws.Range("F8").Value = 0.03168316831683168
$ws.Range("J8").Value = 0.100990099009901
$ws.Range("O8").Value = 0.01188118811881188
$ws.Range("Q8").Value = 0.2455445544554455
$ws.Range("R8").Value = 0.07326732673267326
$ws.Range("S8").Value = 0.4792079207920792

# Row 9
$ws.Range("B9").Value = 0.04590163934426229
$ws.Range("D9").Value = 0.003278688524590164
$ws.Range("E9").Value = 0.003278688524590164
$ws.Range("F9").Value = 0.04262295081967213
$ws.Range("J9").Value = 0.1147540983606557
$ws.Range("O9").Value = 0.01639344262295082
$ws.Range("Q9").Value = 0.2459016393442623
$ws.Range("R9").Value = 0.05245901639344262
$ws.Range("S9").Value = 0.4754098360655737

# Row 10
$ws.Range("B10").Value = 0.07033144704931285
$ws.Range("D10").Value = 0.0169765561843169
$ws.Range("E10").Value = 0.0008084074373484236
$ws.Range("F10").Value = 0.07033144704931285
$ws.Range("J10").Value = 0.1309620048504446
$ws.Range("O10").Value = 0.01212611156022635
$ws.Range("Q10").Value = 0.2506063055780113
$ws.Range("R10").Value = 0.06952303961196443
$ws.Range("S10").Value = 0.3783346806790622

# Row 11
$ws.Range("G11").Value = 0.1434782608695652
$ws.Range("J11").Value = 0.08695652173913043
$ws.Range("K11").Value = 0.1695652173913043
$ws.Range("L11").Value = 0.591304347826087
$ws.Range("S11").Value = 0.008695652173913044

# Row 12
$ws.Range("G12").Value = 0.8085106382978723
$ws.Range("J12").Value = 0.1347517730496454
$ws.Range("K12").Value = 0.007092198581560284
$ws.Range("L12").Value = 0.03546099290780142
$ws.Range("S12").Value = 0.01418439716312057

# Row 15
$ws.Range("H15").Value = 0.1883408071748879
$ws.Range("I15").Value = 0.1165919282511211
$ws.Range("J15").Value = 0.3139013452914798
$ws.Range("K15").Value = 0.07174887892376682
$ws.Range("M15").Value = 0.03139013452914798
$ws.Range("O15").Value = 0.04035874439461883
$ws.Range("S15").Value = 0.2376681614349776

# Row 16
$ws.Range("F16").Value = 0.007936507936507936
$ws.Range("H16").Value = 0.1904761904761905
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.4682539682539683
$ws.Range("K16").Value = 0.06349206349206349
$ws.Range("M16").Value = 0.01587301587301587
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.0873015873015873

# Row 17
$ws.Range("F17").Value = 0.01013513513513514
$ws.Range("H17").Value = 0.1993243243243243
$ws.Range("I17").Value = 0.1402027027027027
$ws.Range("J17").Value = 0.4054054054054054
$ws.Range("K17").Value = 0.06587837837837837
$ws.Range("M17").Value = 0.01858108108108108
$ws.Range("O17").Value = 0.06756756756756757
$ws.Range("S17").Value = 0.0929054054054054

# Row 18
$ws.Range("H18").Value = 0.25
$ws.Range("I18").Value = 0.15
$ws.Range("J18").Value = 0.3375
$ws.Range("K18").Value = 0.0375
$ws.Range("M18").Value = 0.01875
$ws.Range("O18").Value = 0.09375
$ws.Range("S18").Value = 0.1125

# Row 19
$ws.Range("F19").Value = 0.01777434312210201
$ws.Range("H19").Value = 0.2194744976816074
$ws.Range("I19").Value = 0.1244204018547141
$ws.Range("J19").Value = 0.3624420401854714
$ws.Range("K19").Value = 0.09119010819165378
$ws.Range("M19").Value = 0.0170015455950541
$ws.Range("N19").Value = 0.003091190108191654
$ws.Range("O19").Value = 0.06182380216383308
$ws.Range("S19").Value = 0.1027820710973725
